$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Rechercieren wie man Drohne stabiler bekommen" (row 4) -> Fortschritt: Offen -> done
$ws.Range("F4").Value = "done"

# "Simulation fixen (Hindernisse eintragen noch verbuggt)" (row 5) -> Fortschritt: In Arbeit -> done
$ws.Range("F5").Value = "done"

# "Testbilder mit Tiefenkamera" (row 10) -> Fortschritt: offen -> done
$ws.Range("F10").Value = "done"

# "Explodierende Vögel" (row 12) -> Endzeit: (blank) -> "wird wohl leider nie gemacht f"
$ws.Range("E12").Value = "wird wohl leider nie gemacht f"

# "GPS der Drohne auslesen" (row 14) -> Fortschritt: offen -> done
$ws.Range("F14").Value = "done"

# Move the active selection to F16 (matches the saved cursor position in the edited file)
$null = $ws.Range("F16").Select()
